$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the date header row (row 1) from KM1 through LL1, copying
# KM1's number format (style) so the new date cells render like the
# existing ones (numFmtId 14, "m/d/yyyy") without creating new styles.
$ws.Range("KM1").Copy()
$ws.Range("KN1:LL1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 1: new data for columns KN:LL (dates)
$row1 = @(
    44136, 44137, 44138, 44139, 44140, 44141, 44142, 44143, 44144, 44145, 44146, 44147, 44148, 44149, 44150, 44151, 44152, 44153, 44154, 44155, 44156, 44157, 44158, 44159, 44160
)
for ($i = 0; $i -lt 25; $i++) {
    $ws.Cells.Item(1, 300 + $i).Value = $row1[$i]
}

# Row 2: new data for columns KN:LL (mobility values)
$row2 = @(
    53.51, 63.33, 62.67, 61.57, 66.55, 69.89, 65.67, 54.18, 64.599999999999994, 63.26, 61.31, 59.66, 67.5, 61.8, 51.2, 60.21, 60.05, 57.63, 61.45, 64.05, 59.22, 48.94, 62.5, 58.84, 59.18
)
for ($i = 0; $i -lt 25; $i++) {
    $ws.Cells.Item(2, 300 + $i).Value = $row2[$i]
}

# Row 3: new data for columns KN:LL (mobility values)
$row3 = @(
    43.1, 47.29, 43.62, 44.2, 50.18, 58.54, 62.59, 50.02, 49.65, 51.29, 46.52, 46.26, 55.02, 60.93, 43.98, 47.23, 46.72, 46.02, 47.14, 55.9, 57.45, 43.09, 45.59, 46.91, 45.05
)
for ($i = 0; $i -lt 25; $i++) {
    $ws.Cells.Item(3, 300 + $i).Value = $row3[$i]
}

# Row 4: new data for columns KN:LL (mobility values)
$row4 = @(
    52.72, 62.59, 61.18, 59.31, 59.72, 67.989999999999995, 60.38, 49.36, 60.69, 61.92, 63.51, 58.57, 65.27, 59.92, 52.88, 61.1, 62.75, 62.76, 61.12, 66.650000000000006, 61.9, 50.72, 61.7, 64.41, 62.21
)
for ($i = 0; $i -lt 25; $i++) {
    $ws.Cells.Item(4, 300 + $i).Value = $row4[$i]
}

# Row 5: new data for columns KN:LL (mobility values)
$row5 = @(
    33.06, 35.01, 34.32, 34.33, 35.78, 37.869999999999997, 34.67, 31.31, 34.04, 32.43, 34.840000000000003, 30.54, 32.880000000000003, 34.06, 32.57, 35.42, 31.44, 33.909999999999997, 36.35, 36.520000000000003, 31.24, 30.63, 35.08, 35.89, 34.9
)
for ($i = 0; $i -lt 25; $i++) {
    $ws.Cells.Item(5, 300 + $i).Value = $row5[$i]
}
